$wb = $excel.ActiveWorkbook

# The sheets currently named "female" and "male" were mislabeled - swap their
# tab names (this only changes the <sheet name="..."> entries in workbook.xml;
# sheetId/r:id and each sheet's own data stay exactly where they were).
$wsFemale = $wb.Worksheets.Item("female")
$wsMale   = $wb.Worksheets.Item("male")
$wsBoth   = $wb.Worksheets.Item("both")

$wsFemale.Name = "male_tmp__"
$wsMale.Name   = "female"
$wsFemale.Name = "male"

# Every sheet is missing a header label for the state/country name column -
# add it back as "state" in A1.
$wsFemale.Range("A1").Value = "state"
$wsMale.Range("A1").Value   = "state"
$wsBoth.Range("A1").Value   = "state"

# Keep the selection anchored at A1 on the two data sheets ...
$wsMale.Range("A1").Select() | Out-Null
$wsFemale.Range("A1").Select() | Out-Null

# ... and make "both" the active tab, with E3 selected.
$wsBoth.Activate() | Out-Null
$wsBoth.Range("E3").Select() | Out-Null
